# Trade #27 closed at 2026-02-17 15:19:46 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.86
$summary.Range("B4").Value = -0.14
$summary.Range("B5").Value = -0.1
$summary.Range("B6").Value = 27
$summary.Range("B7").Value = 8
$summary.Range("B9").Value = 29.63

# --- Strategy Status sheet ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.86
$status.Range("D4").Value = 27
$status.Range("E4").Value = -0.14
$status.Range("F4").Value = -0.14
$status.Range("G4").Value = 29.63

# --- helper to append the new closed trade row to a trades-log sheet ---
function Add-TradeRow($ws) {
    $ws.Range("A28").Value = 27
    $ws.Range("B28").NumberFormat = "@"
    $ws.Range("B28").Value = "2026-02-17"
    $ws.Range("B28").Style = "Normal"
    $ws.Range("C28").Value = "15:19:39"
    $ws.Range("D28").Value = "MarketMaking"
    $ws.Range("E28").Value = "DOWN"
    $ws.Range("F28").Value = 0.28
    $ws.Range("G28").Value = 0.468193
    $ws.Range("H28").Value = "CLOSED"
    $ws.Range("I28").Value = 67.2119
    $ws.Range("J28").Value = 0.19
    $ws.Range("K28").Value = 99.86
    $ws.Range("L28").Value = 0
    $ws.Range("M28").Value = 0
    $ws.Range("N28").Value = 0.6
    $ws.Range("O28").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P28").Value = "early_exit"
    $ws.Range("Q28").Value = 0.17
}

# --- All Trades sheet ---
$allTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow $allTrades

# --- MarketMaking sheet ---
$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow $marketMaking
